$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Append "*" to the four header labels (Name, Email, Password,
#    Confirm Password) to mark them as required fields.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Name*"
$ws.Range("B1").Value = "Email*"
$ws.Range("C1").Value = "Password*"
$ws.Range("D1").Value = "Confirm Password*"

# ---------------------------------------------------------------------
# 2. The existing header/border style (applied to A1:D2) shrinks its
#    font from 11pt to 8pt and its border colour moves from the old
#    grey (#888888) to a darker grey (#757575).
# ---------------------------------------------------------------------
$existing = $ws.Range("A1:D2")
$existing.Font.Size = 8
$existing.Borders.Color = 7697781  # RGB(117,117,117) -> #757575
$existing.Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 3. Extend formatting across the rest of the used columns (E:U) for
#    rows 1 and 2, mirroring the two border-colour families and the two
#    font families found elsewhere in the sheet. Every range below is a
#    single contiguous block (this COM host does not support multi-area
#    "A1:A1,C1:C1" unions).
# ---------------------------------------------------------------------

# -- Row 1 -----------------------------------------------------------
# E1:Q1 -> Calibri 10pt black text, thin border in the old grey (#888888)
$r1a = $ws.Range("E1:Q1")
$r1a.Font.Name = "Calibri"
$r1a.Font.Size = 10
$r1a.Font.Color = 0
$r1a.Borders.Color = 8947848  # RGB(136,136,136) -> #888888
$r1a.Borders.LineStyle = 1

# R1:U1 -> Calibri 10pt black text, thin border in the new grey (#757575)
$r1b = $ws.Range("R1:U1")
$r1b.Font.Name = "Calibri"
$r1b.Font.Size = 10
$r1b.Font.Color = 0
$r1b.Borders.Color = 7697781  # #757575
$r1b.Borders.LineStyle = 1

# -- Row 2 -------------------------------------------------------------
# Contiguous runs, alternating between the "Calibri black / old grey"
# look and the "theme text / old grey" look, finishing with the
# "theme text / new grey" look for the last four columns.

# E2:F2 -> Calibri 10pt black text, old grey (#888888) border
$g1 = $ws.Range("E2:F2")
$g1.Font.Name = "Calibri"
$g1.Font.Size = 10
$g1.Font.Color = 0
$g1.Borders.Color = 8947848
$g1.Borders.LineStyle = 1

# G2:H2 -> theme text colour 1, Arial 10pt, old grey (#888888) border
$g2 = $ws.Range("G2:H2")
$g2.Font.ThemeColor = 1
$g2.Font.Size = 10
$g2.Borders.Color = 8947848
$g2.Borders.LineStyle = 1

# I2 -> Calibri 10pt black text, old grey (#888888) border
$g3 = $ws.Range("I2:I2")
$g3.Font.Name = "Calibri"
$g3.Font.Size = 10
$g3.Font.Color = 0
$g3.Borders.Color = 8947848
$g3.Borders.LineStyle = 1

# J2 -> theme text colour 1, Arial 10pt, old grey (#888888) border
$g4 = $ws.Range("J2:J2")
$g4.Font.ThemeColor = 1
$g4.Font.Size = 10
$g4.Borders.Color = 8947848
$g4.Borders.LineStyle = 1

# K2:L2 -> Calibri 10pt black text, old grey (#888888) border
$g5 = $ws.Range("K2:L2")
$g5.Font.Name = "Calibri"
$g5.Font.Size = 10
$g5.Font.Color = 0
$g5.Borders.Color = 8947848
$g5.Borders.LineStyle = 1

# M2:P2 -> theme text colour 1, Arial 10pt, old grey (#888888) border
$g6 = $ws.Range("M2:P2")
$g6.Font.ThemeColor = 1
$g6.Font.Size = 10
$g6.Borders.Color = 8947848
$g6.Borders.LineStyle = 1

# Q2 -> Calibri 10pt black text, old grey (#888888) border
$g7 = $ws.Range("Q2:Q2")
$g7.Font.Name = "Calibri"
$g7.Font.Size = 10
$g7.Font.Color = 0
$g7.Borders.Color = 8947848
$g7.Borders.LineStyle = 1

# R2:U2 -> theme text colour 1, Arial 10pt, new grey (#757575) border
$g8 = $ws.Range("R2:U2")
$g8.Font.ThemeColor = 1
$g8.Font.Size = 10
$g8.Borders.Color = 7697781
$g8.Borders.LineStyle = 1

Write-Host "edit complete"
